$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$tf2 = $shape.TextFrame2
$tr2 = $tf2.TextRange
try {
  $para1 = $tr2.Paragraphs(1)
  Write-Host "para1:" $para1
  $para1.Font.Size = 24
} catch {
  Write-Host "fail: $_"
}
